$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write header row (row 1) in the exact order needed so the shared-string
# table is built in the same first-seen order as the target workbook.
$ws.Range("B1").Value = "Owner"
$ws.Range("C1").Value = "Date_start"
$ws.Range("D1").Value = "Date_end"
$ws.Range("E1").Value = "Descr."
$ws.Range("A1").Value = "Topic"

# Topic column (A) for rows 4-6, then owner (B4), then row 7 topic/descr,
# then the remaining owner names - this exact order reproduces the
# original shared-string table order (strings are interned in
# first-seen order as cells are written).
$ws.Range("A4").Value = "QoL"
$ws.Range("A5").Value = "Calls for Service"
$ws.Range("A6").Value = "Stop & Search"
$ws.Range("B4").Value = "Ronak, Niteesh"
$ws.Range("A7").Value = "Issue Analysis & Prioritization"
$ws.Range("E7").Value = "Documentation"
$ws.Range("B7").Value = "Christina"
$ws.Range("B6").Value = "August, Christina"
$ws.Range("B5").Value = "Nikhil, Karthik"

# Date columns C/D for rows 4-7 - set the display format first so the
# engine maps onto the built-in "d-mmm" numFmtId (16) instead of minting a
# custom numFmt, then assign the date values.
$dateRange = $ws.Range("C4:D7")
$dateRange.NumberFormat = "d-mmm"
$ws.Range("C4").Value = [DateTime]"2017-02-10"
$ws.Range("D4").Value = [DateTime]"2017-02-14"
$ws.Range("C5").Value = [DateTime]"2017-02-10"
$ws.Range("D5").Value = [DateTime]"2017-02-14"
$ws.Range("C6").Value = [DateTime]"2017-02-10"
$ws.Range("D6").Value = [DateTime]"2017-02-14"
$ws.Range("C7").Value = [DateTime]"2017-02-10"
$ws.Range("D7").Value = [DateTime]"2017-02-14"

# Column widths (best-fit autosize approximations)
$ws.Columns.Item(1).ColumnWidth = 25.54296875
$ws.Columns.Item(2).ColumnWidth = 14.90625
$ws.Columns.Item(3).ColumnWidth = 9.54296875
$ws.Columns.Item(4).ColumnWidth = 8.90625
$ws.Columns.Item(5).ColumnWidth = 13.81640625

# Selection matches the post-edit active cell in the workbook
[void]$ws.Range("E4").Select()
